$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,5).Value2 = "2026-02-22 23:18:38"
$ws.Cells.Item(3,5).Value2 = "2026-02-22 23:18:41"
$ws.Cells.Item(3,15).Value2 = "4.2 °C"
$ws.Cells.Item(4,5).Value2 = "2026-02-22 23:18:43"
$ws.Cells.Item(4,8).Value2 = "'66%"
$ws.Cells.Item(4,10).Value2 = "1027.1 hPa"
$ws.Cells.Item(4,15).Value2 = "12.0 °C"
$ws.Cells.Item(5,5).Value2 = "2026-02-22 23:18:46"
$ws.Cells.Item(5,8).Value2 = "'30%"
$ws.Cells.Item(5,14).Value2 = "1.7 °C 22:36 TU"
$ws.Cells.Item(5,15).Value2 = "5.8 °C"
$ws.Cells.Item(6,5).Value2 = "2026-02-22 23:18:49"
$ws.Cells.Item(6,15).Value2 = "12.9 °C"
$ws.Cells.Item(7,5).Value2 = "2026-02-22 23:18:51"
$ws.Cells.Item(8,5).Value2 = "2026-02-22 23:18:54"
$ws.Cells.Item(9,5).Value2 = "2026-02-22 23:18:57"
$ws.Cells.Item(10,5).Value2 = "2026-02-22 23:18:59"
$ws.Cells.Item(10,15).Value2 = "9.7 °C"
$ws.Cells.Item(11,5).Value2 = "2026-02-22 23:19:02"
$ws.Cells.Item(12,5).Value2 = "2026-02-22 23:19:04"
$ws.Cells.Item(12,15).Value2 = "9.8 °C"
$ws.Cells.Item(13,5).Value2 = "2026-02-22 23:19:07"
$ws.Cells.Item(13,15).Value2 = "6.4 °C"
$ws.Cells.Item(14,5).Value2 = "2026-02-22 23:19:10"
$ws.Cells.Item(14,15).Value2 = "11.6 °C"
$ws.Cells.Item(15,5).Value2 = "2026-02-22 23:19:12"
$ws.Cells.Item(16,5).Value2 = "2026-02-22 23:19:15"
$ws.Cells.Item(17,5).Value2 = "2026-02-22 23:19:17"
$ws.Cells.Item(17,8).Value2 = "'31%"
$ws.Cells.Item(17,14).Value2 = "7.3 °C 22:55 TU"
$ws.Cells.Item(17,15).Value2 = "9.8 °C"
$ws.Cells.Item(18,5).Value2 = "2026-02-22 23:19:20"
$ws.Cells.Item(18,15).Value2 = "9.8 °C"
$ws.Cells.Item(19,5).Value2 = "2026-02-22 23:19:22"
$ws.Cells.Item(19,15).Value2 = "11.9 °C"
$ws.Cells.Item(20,5).Value2 = "2026-02-22 23:19:25"
$ws.Cells.Item(21,5).Value2 = "2026-02-22 23:19:28"
$ws.Cells.Item(21,8).Value2 = "'60%"
$ws.Cells.Item(21,15).Value2 = "9.1 °C"
$ws.Cells.Item(22,5).Value2 = "2026-02-22 23:19:31"
$ws.Cells.Item(23,5).Value2 = "2026-02-22 23:19:33"
$ws.Cells.Item(23,15).Value2 = "5.7 °C"
$ws.Cells.Item(24,5).Value2 = "2026-02-22 23:19:36"
$ws.Cells.Item(24,8).Value2 = "'85%"
$ws.Cells.Item(24,10).Value2 = "1029.4 hPa"
$ws.Cells.Item(25,5).Value2 = "2026-02-22 23:19:38"
$ws.Cells.Item(26,5).Value2 = "2026-02-22 23:19:41"
$ws.Cells.Item(26,8).Value2 = "'39%"
$ws.Cells.Item(26,15).Value2 = "10.9 °C"
$ws.Cells.Item(27,5).Value2 = "2026-02-22 23:19:44"
$ws.Cells.Item(28,5).Value2 = "2026-02-22 23:19:46"
$ws.Cells.Item(28,8).Value2 = "'66%"
$ws.Cells.Item(28,15).Value2 = "10.2 °C"
$ws.Cells.Item(29,5).Value2 = "2026-02-22 23:19:49"
$ws.Cells.Item(29,15).Value2 = "9.5 °C"
$ws.Cells.Item(30,5).Value2 = "2026-02-22 23:19:52"
$ws.Cells.Item(31,5).Value2 = "2026-02-22 23:19:55"
$ws.Cells.Item(31,8).Value2 = "'60%"
$ws.Cells.Item(31,10).Value2 = "1026.4 hPa"
$ws.Cells.Item(31,12).Value2 = "47.2 km/h - 329º 22:57 TU"
$ws.Cells.Item(32,5).Value2 = "2026-02-22 23:19:57"
$ws.Cells.Item(32,8).Value2 = "'73%"
$ws.Cells.Item(32,15).Value2 = "5.5 °C"
$ws.Cells.Item(33,5).Value2 = "2026-02-22 23:20:00"
$ws.Cells.Item(33,15).Value2 = "8.1 °C"
$ws.Cells.Item(34,5).Value2 = "2026-02-22 23:20:03"
$ws.Cells.Item(35,5).Value2 = "2026-02-22 23:20:06"
$ws.Cells.Item(36,5).Value2 = "2026-02-22 23:20:08"
$ws.Cells.Item(36,8).Value2 = "'78%"
$ws.Cells.Item(36,15).Value2 = "11.5 °C"
$ws.Cells.Item(37,5).Value2 = "2026-02-22 23:20:11"
$ws.Cells.Item(37,15).Value2 = "7.7 °C"
$ws.Cells.Item(38,5).Value2 = "2026-02-22 23:20:14"
$ws.Cells.Item(38,11).Value2 = "15.5 MJ/m2"
$ws.Cells.Item(38,15).Value2 = "11.4 °C"
$ws.Cells.Item(39,5).Value2 = "2026-02-22 23:20:16"
$ws.Cells.Item(40,5).Value2 = "2026-02-22 23:20:19"
$ws.Cells.Item(40,8).Value2 = "'59%"
$ws.Cells.Item(40,15).Value2 = "9.5 °C"
$ws.Cells.Item(41,5).Value2 = "2026-02-22 23:20:21"
$ws.Cells.Item(42,5).Value2 = "2026-02-22 23:20:24"
$ws.Cells.Item(42,15).Value2 = "10.4 °C"
$ws.Cells.Item(43,5).Value2 = "2026-02-22 23:20:26"
$ws.Cells.Item(44,5).Value2 = "2026-02-22 23:20:29"
$ws.Cells.Item(44,15).Value2 = "2.7 °C"
$ws.Cells.Item(45,5).Value2 = "2026-02-22 23:20:31"
$ws.Cells.Item(45,8).Value2 = "'56%"
$ws.Cells.Item(45,10).Value2 = "1029.2 hPa"
$ws.Cells.Item(45,15).Value2 = "8.4 °C"
$ws.Cells.Item(46,5).Value2 = "2026-02-22 23:20:34"
$ws.Cells.Item(46,15).Value2 = "8.8 °C"
